$d = $word.ActiveDocument

$pairs = @(
    @("333÷9=37, 0", "419÷2=209, 1"),
    @("379÷8=47, 3", "574÷5=114, 4"),
    @("997÷3=332, 1", "783÷9=87, 0"),
    @("967÷5=193, 2", "878÷6=146, 2"),
    @("726÷5=145, 1", "637÷2=318, 1"),
    @("850÷5=170, 0", "677÷4=169, 1"),
    @("647÷5=129, 2", "613÷3=204, 1"),
    @("445÷3=148, 1", "500÷6=83, 2"),
    @("751÷2=375, 1", "167÷9=18, 5"),
    @("941÷8=117, 5", "337÷9=37, 4"),
    @("455÷2=227, 1", "359÷3=119, 2"),
    @("493÷8=61, 5", "503÷5=100, 3"),
    @("562÷9=62, 4", "442÷7=63, 1"),
    @("379÷9=42, 1", "252÷7=36, 0"),
    @("866÷2=433, 0", "784÷3=261, 1"),
    @("453÷3=151, 0", "192÷6=32, 0"),
    @("520÷4=130, 0", "117÷9=13, 0"),
    @("986÷4=246, 2", "702÷9=78, 0"),
    @("876÷5=175, 1", "257÷2=128, 1"),
    @("177÷9=19, 6", "856÷2=428, 0"),
    @("595÷3=198, 1", "654÷7=93, 3"),
    @("228÷4=57, 0", "266÷8=33, 2"),
    @("405÷6=67, 3", "480÷3=160, 0"),
    @("753÷9=83, 6", "320÷6=53, 2"),
    @("924÷9=102, 6", "257÷8=32, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
